$d = $word.ActiveDocument

# Locate the empty paragraph that precedes the insertion point.
# It is the paragraph right before the final (last) paragraph in the document body.
$anchorIndex = $d.Paragraphs.Count - 1
$anchor = $d.Paragraphs.Item($anchorIndex)

# --- New paragraph 1 ---
$anchor.Range.InsertParagraphAfter() | Out-Null
$anchorIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($anchorIndex)
$xml1 = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$newPara.Range.InsertXML($xml1) | Out-Null
$anchor = $newPara

# --- New paragraph 2 ---
$anchor.Range.InsertParagraphAfter() | Out-Null
$anchorIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($anchorIndex)
$xml2 = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="720"/>
          <w:tab w:val="left" w:pos="1440"/>
          <w:tab w:val="left" w:pos="2160"/>
          <w:tab w:val="left" w:pos="2880"/>
          <w:tab w:val="left" w:pos="3600"/>
          <w:tab w:val="left" w:pos="4320"/>
          <w:tab w:val="left" w:pos="5040"/>
          <w:tab w:val="left" w:pos="5760"/>
          <w:tab w:val="left" w:pos="6480"/>
          <w:tab w:val="left" w:pos="7200"/>
          <w:tab w:val="left" w:pos="7920"/>
          <w:tab w:val="left" w:pos="8640"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">Optogenetics seem </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">to be such a paradigm shift </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">in biomedical engineering. Your comment about being invasive is probably a major roadblock to transition easily the therapies to humans. </w:t>
      </w:r>
    </w:p>
'@
$newPara.Range.InsertXML($xml2) | Out-Null
$anchor = $newPara

# --- New paragraph 3 ---
$anchor.Range.InsertParagraphAfter() | Out-Null
$anchorIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($anchorIndex)
$xml3 = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">Still when you </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>mention</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> restoring motor functions from paralysis, it seems that the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>stakes</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>are</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> so </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>important that</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>t</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">here </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>will be a</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">strong </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>motivation to make it successful</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$newPara.Range.InsertXML($xml3) | Out-Null
$anchor = $newPara

# --- New paragraph 4 ---
$anchor.Range.InsertParagraphAfter() | Out-Null
$anchorIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($anchorIndex)
$xml4 = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$newPara.Range.InsertXML($xml4) | Out-Null
$anchor = $newPara

# --- New paragraph 5 ---
$anchor.Range.InsertParagraphAfter() | Out-Null
$anchorIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($anchorIndex)
$xml5 = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">Reading your post about using </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>optogentic</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> therapies to reduce pain made me think that it could be useful for some surgeries where you want the patient to be </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>conscious</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> without using anesthesia or other surgeries </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">where even administering the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>anesthesia</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> like with an epidural, could be painful. The ethical considerations about editing genes though, I do not think it is specific to optogenetics. However, I could imagine cases for which decreasing levels of pain could be used outside of medical contexts. Societies will have to address the legal implications but it might be one of th</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>ese</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t>situations</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve"> where </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">pace of </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Helvetica" w:cs="Helvetica"/>
        </w:rPr>
        <w:t xml:space="preserve">progress is too fast for them to handle quickly enough their implications (See AI fake news, or responsible AI and other related topics). </w:t>
      </w:r>
    </w:p>
'@
$newPara.Range.InsertXML($xml5) | Out-Null
$anchor = $newPara

Write-Output "Inserted paragraphs. Total paragraphs now:"
Write-Output $d.Paragraphs.Count
